$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the existing row 471 (the 2020-12-02 "Frutilla" block),
# shifting all rows from 471 down through 481 to 474 through 484.
$ws.Range("A471:T473").EntireRow.Insert()

# Populate the 3 newly inserted rows with the new weekly report (Fecha = 44448,
# i.e. 2021-09-09) for "Frutilla" quality grades Primera / Segunda / Tercera.
$newRows = @(
    @{ Row = 471; Quality = "Primera"; Vol = 65;  Min = 20000; Max = 22000; Avg = 21077; PerKg = 3011 },
    @{ Row = 472; Quality = "Segunda"; Vol = 50;  Min = 15000; Max = 16000; Avg = 15600; PerKg = 2229 },
    @{ Row = 473; Quality = "Tercera"; Vol = 50;  Min = 8000;  Max = 10000; Avg = 9000;  PerKg = 1286 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 9
    $ws.Cells.Item($row, 2).Value = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($row, 3).Value = "Metropolitana"
    $ws.Cells.Item($row, 4).Value = 44448
    $ws.Cells.Item($row, 5).Value = 13
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100101
    $ws.Cells.Item($row, 8).Value = "Berries"
    $ws.Cells.Item($row, 9).Value = 100112025
    $ws.Cells.Item($row, 10).Value = "Frutilla"
    $ws.Cells.Item($row, 11).Value = "Sin especificar"
    $ws.Cells.Item($row, 12).Value = $r.Quality
    $ws.Cells.Item($row, 13).Value = $r.Vol
    $ws.Cells.Item($row, 14).Value = $r.Min
    $ws.Cells.Item($row, 15).Value = $r.Max
    $ws.Cells.Item($row, 16).Value = $r.Avg
    $ws.Cells.Item($row, 17).Value = "`$/bandeja 7 kilos"
    $ws.Cells.Item($row, 18).Value = "Provincia de Melipilla"
    $ws.Cells.Item($row, 19).Value = $r.PerKg
    $ws.Cells.Item($row, 20).Value = 7
}

# The report dated 44400 (2021-07-23) that lands on row 481 after the shift was
# actually logged a day earlier (44399 / 2021-07-22); correct its Fecha.
$ws.Cells.Item(481, 4).Value = 44399
